# connected location database with from output
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Location database lookup pulled a MaterialNo of 999 into row 12 (A12)
$ws.Range("A12").Value = 999

# Quantity from output pulled 9999999 into row 15 (B15)
$ws.Range("B15").Value = 9999999

# Selection moved up one row (B8 -> B7) as the cursor landed on the row
# that now has data
$ws.Range("B7").Select()
